$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# --- Step 1: bump the printed date in A1 by one day (45308 -> 45309) ---
$ws.Range("A1").Value = 45309

# --- Step 2: update the three price cells ---
$ws.Range("D33").Value = 1305
$ws.Range("D34").Value = 1275.478
$ws.Range("D35").Value = 949.728

